$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# "Alter bei AN-Anteil PV" (age relevant for the employer's long-term-care
# insurance contribution): insert a new row above the old row 41
# ("wohnhaft Sachsen?"), pushing rows 41:48 down to 42:49, and fill it
# with the new question + default answer.
$ws.Rows.Item(41).Insert()

$ws.Range("A41").Value = "juenger als 23 oder vor 1940 geboren?"
$ws.Range("B41").Value = "nein"

# B41 picks up the same "ja/nein" number format used by the row above it
# (row 40, B40) rather than the default General format.
$ws.Range("B41").NumberFormat = $ws.Range("B40").NumberFormat

# Restore the view: scrolled down so row 19 is at the top, with A33
# selected as the active cell.
try {
    $excel.ActiveWindow.ScrollRow = 19
} catch {}
$ws.Range("A33").Select() | Out-Null
